# 9.2.1.xlsx: add a new "2021" column (R) that mirrors the existing
# "2020" column (Q) formatting, with the new year's data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (thin spacer/border row above the header) - copy Q2's formatting,
# no value, into the new R2 cell.
$ws.Range("Q2").Copy($ws.Range("R2"))

# Row 3 (year header row) - copy Q3's formatting into R3, then set 2021.
$ws.Range("Q3").Copy($ws.Range("R3"))
$ws.Range("R3").Value = 2021

# Row 4 (GVA share of manufacturing output in GDP, %) - copy Q4's
# formatting into R4, then set the 2021 figure.
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 13.5

# Row 5 (GVA of manufacturing industry in GDP per capita) - copy Q5's
# formatting into R5, then set the 2021 figure.
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 15.1

# Match the author's final cursor position/selection.
$ws.Range("T3").Select()
